$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.703.64"
$ws.Range("E2").Value = "  +0.67%  "

$ws.Range("D3").Value = "'1.961.73"
$ws.Range("E3").Value = "  +1.01%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'244.98"
$ws.Range("E5").Value = "  +1.01%  "

$ws.Range("D6").Value = "'0.618"
$ws.Range("E6").Value = "  +2.65%  "

$ws.Range("D7").Value = "'61.66"
$ws.Range("E7").Value = "  +7.98%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D10").Value = "'0.0798"
$ws.Range("E10").Value = "  -6.51%  "

$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("D12").Value = "'14.27"
$ws.Range("E12").Value = "  +6.07%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.839"
$ws.Range("E13").Value = "  +3.88%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'21.98"
$ws.Range("E14").Value = "  +2.48%  "

$ws.Range("D15").Value = "'2.236.17"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("D16").Value = "'5.31"
$ws.Range("E16").Value = "  +3.25%  "

$ws.Range("D17").Value = "'1.960.92"
$ws.Range("E17").Value = "  +1.35%  "

$ws.Range("D18").Value = "'36.642.22"
$ws.Range("E18").Value = "  +0.69%  "

$ws.Range("D19").Value = "'70.03"
$ws.Range("E19").Value = "  +1.34%  "

$ws.Range("D20").Value = "0.0₃0855"
$ws.Range("E20").Value = "  -0.91%  "

$ws.Range("D21").Value = "'230.56"
$ws.Range("E21").Value = "  +1.51%  "

$ws.Range("D22").Value = "'5.09"
$ws.Range("E22").Value = "  +1.99%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("E24").Value = "  +5.89%  "

$ws.Range("E25").Value = "  +3.17%  "

$ws.Range("D26").Value = "'0.142"
$ws.Range("E26").Value = "  +5.86%  "

$ws.Range("D27").Value = "'9.23"
$ws.Range("E27").Value = "  +0.90%  "

$ws.Range("D28").Value = "'160.70"
$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("D29").Value = "'19.47"
$ws.Range("E29").Value = "  +1.04%  "

$ws.Range("D30").Value = "'1.31"
$ws.Range("E30").Value = "  +18.37%  "

$ws.Range("E31").Value = "  +1.60%  "

$ws.Range("D32").Value = "'4.79"
$ws.Range("E32").Value = "  +5.17%  "

$ws.Range("D33").Value = "'0.0620"
$ws.Range("E33").Value = "  -0.15%  "

$ws.Range("E34").Value = "  +7.48%  "

$ws.Range("D35").Value = "'3.55"
$ws.Range("E35").Value = "  +16.17%  "

$ws.Range("D36").Value = "'2.27"
$ws.Range("E36").Value = "  +5.16%  "

$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").Value = "'1.78"
$ws.Range("E38").Value = "  -1.06%  "

$ws.Range("E39").Value = "  -7.38%  "

$ws.Range("D40").Value = "'0.0986"
$ws.Range("E40").Value = "  -0.38%  "

$ws.Range("E41").Value = "  +0.91%  "

$ws.Range("E42").Value = "  +2.54%  "

$ws.Range("E43").Value = "  +1.35%  "

$ws.Range("E44").Value = "  +3.69%  "

$ws.Range("D45").Value = "'1.372.02"
$ws.Range("E45").Value = "  +2.28%  "

$ws.Range("D46").Value = "'88.97"
$ws.Range("E46").Value = "  +3.14%  "

$ws.Range("E47").Value = "  +1.83%  "

$ws.Range("D48").Value = "'7.18"
$ws.Range("E48").Value = "  +1.08%  "

$ws.Range("E49").Value = "  +0.47%  "

$ws.Range("D50").Value = "'45.55"
$ws.Range("E50").Value = "  +4.78%  "

$ws.Range("D51").Value = "'2.130.89"
$ws.Range("E51").Value = "  +0.77%  "
